# Add New York state hospitalization data for the new reporting day
# (serial date 43929) and backfill the ICU/Intubations figures that
# came in late for the previous day (serial date 43928, row 24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction to an existing value (C23: 89 -> 35)
$ws.Range("C23").Value = 35

# Late-arriving ICU/Intubations counts for the existing last row (row 24)
$ws.Range("C24").Value = 302
$ws.Range("D24").Value = 94

# New row of data: copy the date formatting from the row above, then
# fill in the new day's figures.
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)

$ws.Range("A25").Value = 43929
$ws.Range("B25").Value = 200
$ws.Range("C25").Value = 84
$ws.Range("D25").Value = 88

# Update the active selection to match the new last-populated cell.
[void]$ws.Range("E24").Select()

# Widen the sheet-tabs area of the window (cosmetic window setting).
$excel.ActiveWindow.TabRatio = 995
